$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the stray row 40 cell that held "Error en activacion"
#    (its text becomes the new "Nit" header, see step 4 below).
$ws.Rows.Item(40).Delete()

# 2. Column A ("Iccid") previously carried a text-number-format column
#    style (style index 3 in the original styles.xml). Drop it now that
#    the used range has shrunk back to row 1 only, so this only touches
#    the header row.
$ws.Columns.Item(1).ClearFormats()

# 3. Shift the "Mensaje" header from C1 to D1, carrying its formatting
#    (bold font + border + center/top alignment) along with it.
#    -4122 = xlPasteFormats
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D1").Value = $ws.Range("C1").Value2

# 4. Put the new "Nit" header into C1, reusing its existing header style.
$ws.Range("C1").Value = "Nit"

# 5. Re-apply the standard header style (bold + border) to A1, since it
#    lost its style when the column-level format was cleared in step 2.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A1").Value = "Iccid"

# 6. Column B now needs a bestFit-style custom width.
$ws.Columns.Item(2).ColumnWidth = 10.2

$excel.CutCopyMode = $false
